# Generate Report for Handback
# Row 7 on both the "zh-cn" and "de-de" sheets corresponds to the
# 34a191ca-927e-4061-9901-a63b7b78d9e9 source file. A handback was
# received for it, but its commit is stale relative to the latest
# handoff, so the report now fills in the "Latest Target File",
# "Latest Handback File", "Latest Handback DateTime" and "Error Detail"
# columns for that row (I, J, K, P) on each language sheet, and widens
# the "Error Detail" column to fit the message.

$wb = $excel.ActiveWorkbook

$mdTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f0f7444fa854c6be1be09f26985d406ec2b8a992/e2e/34a191ca-927e-4061-9901-a63b7b78d9e9.md"
$mdDisplay = "34a191ca-927e-4061-9901-a63b7b78d9e9.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f79ec6c3efede5eb3a1ca2fe91ae0749be90ddc5/e2e/34a191ca-927e-4061-9901-a63b7b78d9e9.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f0f7444fa854c6be1be09f26985d406ec2b8a992/e2e/34a191ca-927e-4061-9901-a63b7b78d9e9.md."

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("J7").Value = "34a191ca-927e-4061-9901-a63b7b78d9e9.1f70c007d909377dfe9dd91f6bbcef419dc48419.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-29 02:43:11"
$wsZh.Range("P7").Value = $errorDetail

$wsZh.Range("I7").Value = $mdDisplay
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $mdTarget, "", "", $mdDisplay) | Out-Null

$wsZh.Columns.Item(16).ColumnWidth = 39.083333333333336

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("J7").Value = "34a191ca-927e-4061-9901-a63b7b78d9e9.1f70c007d909377dfe9dd91f6bbcef419dc48419.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-29 02:43:18"
$wsDe.Range("P7").Value = $errorDetail

$wsDe.Range("I7").Value = $mdDisplay
$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $mdTarget, "", "", $mdDisplay) | Out-Null

$wsDe.Columns.Item(16).ColumnWidth = 39.083333333333336

Write-Host "Handback report generated for 34a191ca-927e-4061-9901-a63b7b78d9e9 (zh-cn, de-de)"
